$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F3").Value = 0
$ws.Range("F7").Value = 0
$ws.Range("F8").Value = -6
$ws.Range("F12").Value = -1
$ws.Range("F13").Value = -5
$ws.Range("F17").Value = -1
$ws.Range("F19").Value = -6
$ws.Range("F21").Value = -4
$ws.Range("F22").Value = -9
$ws.Range("F23").Value = -6
$ws.Range("F24").Value = -9
$ws.Range("F26").Value = -5
$ws.Range("F29").Value = 1
$ws.Range("F30").Value = -5
$ws.Range("F34").Value = -3
